$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 633.8  # H9: 621.4286 -> 633.8
$ws.Cells.Item(9, 9).Value = 684.875  # I9: 650.3333 -> 684.875
$ws.Cells.Item(9, 10).Value = 429.5  # J9: 448 -> 429.5
$ws.Cells.Item(9, 11).Value = 684.875  # K9: 650.3333 -> 684.875
$ws.Cells.Item(9, 12).Value = 429.5  # L9: 448 -> 429.5
$ws.Cells.Item(9, 13).Value = -515.875  # M9: -481.3333 -> -515.875
$ws.Cells.Item(9, 14).Value = -767.5  # N9: -786 -> -767.5
$ws.Cells.Item(17, 8).Value = 1548  # H17: 1580.1538 -> 1548
$ws.Cells.Item(17, 10).Value = 1548  # J17: 1580.1538 -> 1548
$ws.Cells.Item(17, 12).Value = 4644  # L17: 4740.4614 -> 4644
$ws.Cells.Item(17, 14).Value = -4980  # N17: -5076.4614 -> -4980
$ws.Cells.Item(19, 8).Value = 3258.5  # H19: 3297.0908 -> 3258.5
$ws.Cells.Item(19, 9).Value = 2381.3333  # I19: 2378.6667 -> 2381.3333
$ws.Cells.Item(19, 10).Value = 4574.25  # J19: 4399.2 -> 4574.25
$ws.Cells.Item(19, 11).Value = 2381.3333  # K19: 2378.6667 -> 2381.3333
$ws.Cells.Item(19, 12).Value = 4574.25  # L19: 4399.2 -> 4574.25
$ws.Cells.Item(19, 13).Value = -2206.3333  # M19: -2203.6667 -> -2206.3333
$ws.Cells.Item(19, 14).Value = -4924.25  # N19: -4749.2 -> -4924.25
$ws.Cells.Item(38, 8).Value = 324.375  # H38: 351.85715 -> 324.375
$ws.Cells.Item(38, 9).Value = 324.375  # I38: 351.85715 -> 324.375
$ws.Cells.Item(38, 11).Value = 973.125  # K38: 1055.57145 -> 973.125
$ws.Cells.Item(38, 13).Value = -601.125  # M38: -683.5714499999999 -> -601.125
$ws.Cells.Item(80, 8).Value = 3029  # H80: 2968.8 -> 3029
$ws.Cells.Item(80, 9).Value = 5524.5  # I80: 6866.3335 -> 5524.5
$ws.Cells.Item(80, 10).Value = 1365.3334  # J80: 1298.4286 -> 1365.3334
$ws.Cells.Item(80, 11).Value = 16573.5  # K80: 20599.0005 -> 16573.5
$ws.Cells.Item(80, 12).Value = 4096.0002  # L80: 3895.2858 -> 4096.0002
$ws.Cells.Item(80, 13).Value = -15575.5  # M80: -19601.0005 -> -15575.5
$ws.Cells.Item(80, 14).Value = -6092.0002  # N80: -5891.2858 -> -6092.0002
$ws.Cells.Item(83, 8).Value = 3029  # H83: 2968.8 -> 3029
$ws.Cells.Item(83, 9).Value = 5524.5  # I83: 6866.3335 -> 5524.5
$ws.Cells.Item(83, 10).Value = 1365.3334  # J83: 1298.4286 -> 1365.3334
$ws.Cells.Item(83, 11).Value = 49720.5  # K83: 61797.0015 -> 49720.5
$ws.Cells.Item(83, 12).Value = 12288.0006  # L83: 11685.8574 -> 12288.0006
$ws.Cells.Item(83, 13).Value = -44728.5  # M83: -56805.0015 -> -44728.5
$ws.Cells.Item(83, 14).Value = -22272.0006  # N83: -21669.8574 -> -22272.0006
$ws.Cells.Item(88, 8).Value = 415291.16  # H88: 361448.9 -> 415291.16
$ws.Cells.Item(88, 9).Value = 17664.25  # I88: 26096.75 -> 17664.25
$ws.Cells.Item(88, 10).Value = 1011731.5  # J88: 540303.4 -> 1011731.5
$ws.Cells.Item(88, 11).Value = 17664.25  # K88: 26096.75 -> 17664.25
$ws.Cells.Item(88, 12).Value = 1011731.5  # L88: 540303.4 -> 1011731.5
$ws.Cells.Item(88, 13).Value = -17258.25  # M88: -25690.75 -> -17258.25
$ws.Cells.Item(88, 14).Value = -1012543.5  # N88: -541115.4 -> -1012543.5
$ws.Cells.Item(91, 8).Value = 415291.16  # H91: 361448.9 -> 415291.16
$ws.Cells.Item(91, 9).Value = 17664.25  # I91: 26096.75 -> 17664.25
$ws.Cells.Item(91, 10).Value = 1011731.5  # J91: 540303.4 -> 1011731.5
$ws.Cells.Item(91, 11).Value = 17664.25  # K91: 26096.75 -> 17664.25
$ws.Cells.Item(91, 12).Value = 1011731.5  # L91: 540303.4 -> 1011731.5
$ws.Cells.Item(91, 13).Value = -16260.25  # M91: -24692.75 -> -16260.25
$ws.Cells.Item(91, 14).Value = -1014539.5  # N91: -543111.4 -> -1014539.5
$ws.Cells.Item(105, 8).Value = 52223  # H105: 36667.8 -> 52223
$ws.Cells.Item(105, 10).Value = 56667.8  # J105: 38335 -> 56667.8
$ws.Cells.Item(105, 12).Value = 56667.8  # L105: 38335 -> 56667.8
$ws.Cells.Item(105, 14).Value = -63655.8  # N105: -45323 -> -63655.8
$ws.Cells.Item(137, 8).Value = 422911.28  # H137: 461270.53 -> 422911.28
$ws.Cells.Item(137, 9).Value = 558879.0600000001  # I137: 558960.2 -> 558879.0600000001
$ws.Cells.Item(137, 10).Value = 15008  # J137: 21667.25 -> 15008
$ws.Cells.Item(137, 11).Value = 1676637.18  # K137: 1676880.6 -> 1676637.18
$ws.Cells.Item(137, 12).Value = 45024  # L137: 65001.75 -> 45024
$ws.Cells.Item(137, 13).Value = -1674087.18  # M137: -1674330.6 -> -1674087.18
$ws.Cells.Item(137, 14).Value = -50124  # N137: -70101.75 -> -50124
$ws.Cells.Item(138, 8).Value = 4228.1685  # H138: 4216.5054 -> 4228.1685
$ws.Cells.Item(138, 10).Value = 4271.0234  # J138: 4257.989 -> 4271.0234
$ws.Cells.Item(138, 12).Value = 12813.0702  # L138: 12773.967 -> 12813.0702
$ws.Cells.Item(138, 14).Value = -23093.0702  # N138: -23053.967 -> -23093.0702

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 52498.133  # H32: 56462.105 -> 52498.133
$ws.Cells.Item(32, 9).Value = 50219.223  # I32: 54233.547 -> 50219.223
$ws.Cells.Item(32, 10).Value = 57625.688  # J32: 61715.145 -> 57625.688
$ws.Cells.Item(32, 11).Value = 50219.223  # K32: 54233.547 -> 50219.223
$ws.Cells.Item(32, 12).Value = 57625.688  # L32: 61715.145 -> 57625.688
$ws.Cells.Item(32, 13).Value = -49932.223  # M32: -53946.547 -> -49932.223
$ws.Cells.Item(32, 14).Value = -58199.688  # N32: -62289.145 -> -58199.688
$ws.Cells.Item(45, 8).Value = 20835380  # H45: 20001930 -> 20835380
$ws.Cells.Item(45, 9).Value = 23811460  # I45: 21740990 -> 23811460
$ws.Cells.Item(45, 10).Value = 2823.3333  # J45: 2735 -> 2823.3333
$ws.Cells.Item(45, 11).Value = 23811460  # K45: 21740990 -> 23811460
$ws.Cells.Item(45, 12).Value = 2823.3333  # L45: 2735 -> 2823.3333
$ws.Cells.Item(45, 13).Value = -23811083  # M45: -21740613 -> -23811083
$ws.Cells.Item(45, 14).Value = -3577.3333  # N45: -3489 -> -3577.3333
$ws.Cells.Item(61, 8).Value = 10875983  # H61: 11911503 -> 10875983
$ws.Cells.Item(61, 9).Value = 4872.4  # I61: 5024.0713 -> 4872.4
$ws.Cells.Item(61, 10).Value = 31259316  # J61: 35724460 -> 31259316
$ws.Cells.Item(61, 11).Value = 4872.4  # K61: 5024.0713 -> 4872.4
$ws.Cells.Item(61, 12).Value = 31259316  # L61: 35724460 -> 31259316
$ws.Cells.Item(61, 13).Value = -4660.4  # M61: -4812.0713 -> -4660.4
$ws.Cells.Item(61, 14).Value = -31259740  # N61: -35724884 -> -31259740
$ws.Cells.Item(74, 8).Value = 8149.7812  # H74: 7536.029 -> 8149.7812
$ws.Cells.Item(74, 9).Value = 2721.739  # I74: 2521.8462 -> 2721.739
$ws.Cells.Item(74, 11).Value = 2721.739  # K74: 2521.8462 -> 2721.739
$ws.Cells.Item(74, 13).Value = -1847.739  # M74: -1647.8462 -> -1847.739
$ws.Cells.Item(77, 8).Value = 8149.7812  # H77: 7536.029 -> 8149.7812
$ws.Cells.Item(77, 9).Value = 2721.739  # I77: 2521.8462 -> 2721.739
$ws.Cells.Item(77, 11).Value = 13608.695  # K77: 12609.231 -> 13608.695
$ws.Cells.Item(77, 13).Value = -9240.695  # M77: -8241.231 -> -9240.695
$ws.Cells.Item(136, 8).Value = 10875983  # H136: 11911503 -> 10875983
$ws.Cells.Item(136, 9).Value = 4872.4  # I136: 5024.0713 -> 4872.4
$ws.Cells.Item(136, 10).Value = 31259316  # J136: 35724460 -> 31259316
$ws.Cells.Item(136, 11).Value = 14617.2  # K136: 15072.2139 -> 14617.2
$ws.Cells.Item(136, 12).Value = 93777948  # L136: 107173380 -> 93777948
$ws.Cells.Item(136, 13).Value = -12067.2  # M136: -12522.2139 -> -12067.2
$ws.Cells.Item(136, 14).Value = -93783048  # N136: -107178480 -> -93783048

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 5210948.5  # H134: 5558310.5 -> 5210948.5
$ws.Cells.Item(134, 9).Value = 1394.9584  # I134: 1442.6522 -> 1394.9584
$ws.Cells.Item(134, 10).Value = 20839608  # J134: 23816590 -> 20839608
$ws.Cells.Item(134, 11).Value = 4184.8752  # K134: 4327.9566 -> 4184.8752
$ws.Cells.Item(134, 12).Value = 62518824  # L134: 71449770 -> 62518824
$ws.Cells.Item(134, 13).Value = -1649.8752  # M134: -1792.9566 -> -1649.8752
$ws.Cells.Item(134, 14).Value = -62523894  # N134: -71454840 -> -62523894

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8011.4287  # H31: 8020.952 -> 8011.4287
$ws.Cells.Item(31, 9).Value = 4243.3  # I31: 4212.091 -> 4243.3
$ws.Cells.Item(31, 10).Value = 11437  # J31: 12210.7 -> 11437
$ws.Cells.Item(31, 11).Value = 4243.3  # K31: 4212.091 -> 4243.3
$ws.Cells.Item(31, 12).Value = 11437  # L31: 12210.7 -> 11437
$ws.Cells.Item(31, 13).Value = -3948.3  # M31: -3917.091 -> -3948.3
$ws.Cells.Item(31, 14).Value = -12027  # N31: -12800.7 -> -12027
$ws.Cells.Item(34, 8).Value = 8011.4287  # H34: 8020.952 -> 8011.4287
$ws.Cells.Item(34, 9).Value = 4243.3  # I34: 4212.091 -> 4243.3
$ws.Cells.Item(34, 10).Value = 11437  # J34: 12210.7 -> 11437
$ws.Cells.Item(34, 11).Value = 4243.3  # K34: 4212.091 -> 4243.3
$ws.Cells.Item(34, 12).Value = 11437  # L34: 12210.7 -> 11437
$ws.Cells.Item(34, 13).Value = -4041.3  # M34: -4010.091 -> -4041.3
$ws.Cells.Item(34, 14).Value = -11841  # N34: -12614.7 -> -11841
$ws.Cells.Item(35, 8).Value = 981.6667  # H35: 1508.3334 -> 981.6667
$ws.Cells.Item(35, 9).Value = 981.6667  # I35: 1508.3334 -> 981.6667
$ws.Cells.Item(35, 11).Value = 981.6667  # K35: 1508.3334 -> 981.6667
$ws.Cells.Item(35, 13).Value = -687.6667  # M35: -1214.3334 -> -687.6667
$ws.Cells.Item(87, 8).Value = 110000  # H87: 55000 -> 110000
$ws.Cells.Item(87, 10).Value = 110000  # J87: 55000 -> 110000
$ws.Cells.Item(87, 12).Value = 110000  # L87: 55000 -> 110000
$ws.Cells.Item(87, 14).Value = -112372  # N87: -57372 -> -112372
$ws.Cells.Item(90, 8).Value = 110000  # H90: 55000 -> 110000
$ws.Cells.Item(90, 10).Value = 110000  # J90: 55000 -> 110000
$ws.Cells.Item(90, 12).Value = 330000  # L90: 165000 -> 330000
$ws.Cells.Item(90, 14).Value = -341856  # N90: -176856 -> -341856
$ws.Cells.Item(94, 8).Value = 4194.3794  # H94: 4412.926 -> 4194.3794
$ws.Cells.Item(94, 9).Value = 3996  # I94: 4137.643 -> 3996
$ws.Cells.Item(94, 10).Value = 4406.9287  # J94: 4709.385 -> 4406.9287
$ws.Cells.Item(94, 11).Value = 3996  # K94: 4137.643 -> 3996
$ws.Cells.Item(94, 12).Value = 4406.9287  # L94: 4709.385 -> 4406.9287
$ws.Cells.Item(94, 13).Value = -3545  # M94: -3686.643 -> -3545
$ws.Cells.Item(94, 14).Value = -5308.9287  # N94: -5611.385 -> -5308.9287

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 921.6667  # H34: 1079.5 -> 921.6667
$ws.Cells.Item(34, 10).Value = 2000  # J34: 2250 -> 2000
$ws.Cells.Item(34, 12).Value = 6000  # L34: 6750 -> 6000
$ws.Cells.Item(34, 14).Value = -6168  # N34: -6918 -> -6168
$ws.Cells.Item(55, 8).Value = 6999.3335  # H55: 5199.6 -> 6999.3335
$ws.Cells.Item(55, 10).Value = 0  # J55: 2500 -> 0
$ws.Cells.Item(55, 12).Value = 0  # L55: 7500 -> 0
$ws.Cells.Item(55, 14).ClearContents()  # N55: -7854 -> (cleared)
$ws.Cells.Item(131, 8).Value = 14335.333  # H131: 10667.857 -> 14335.333
$ws.Cells.Item(131, 9).Value = 9455  # I131: 11166 -> 9455
$ws.Cells.Item(131, 10).Value = 17588.889  # J131: 10512.1875 -> 17588.889
$ws.Cells.Item(131, 11).Value = 28365  # K131: 33498 -> 28365
$ws.Cells.Item(131, 12).Value = 52766.667  # L131: 31536.5625 -> 52766.667
$ws.Cells.Item(131, 13).Value = -23325  # M131: -28458 -> -23325
$ws.Cells.Item(131, 14).Value = -62846.667  # N131: -41616.5625 -> -62846.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5550.769  # H70: 5364.9644 -> 5550.769
$ws.Cells.Item(70, 9).Value = 5878.1763  # I70: 5712.6665 -> 5878.1763
$ws.Cells.Item(70, 10).Value = 4932.3335  # J70: 4739.1 -> 4932.3335
$ws.Cells.Item(70, 11).Value = 5878.1763  # K70: 5712.6665 -> 5878.1763
$ws.Cells.Item(70, 12).Value = 4932.3335  # L70: 4739.1 -> 4932.3335
$ws.Cells.Item(70, 13).Value = -5608.1763  # M70: -5442.6665 -> -5608.1763
$ws.Cells.Item(70, 14).Value = -5472.3335  # N70: -5279.1 -> -5472.3335
$ws.Cells.Item(73, 8).Value = 5550.769  # H73: 5364.9644 -> 5550.769
$ws.Cells.Item(73, 9).Value = 5878.1763  # I73: 5712.6665 -> 5878.1763
$ws.Cells.Item(73, 10).Value = 4932.3335  # J73: 4739.1 -> 4932.3335
$ws.Cells.Item(73, 11).Value = 5878.1763  # K73: 5712.6665 -> 5878.1763
$ws.Cells.Item(73, 12).Value = 4932.3335  # L73: 4739.1 -> 4932.3335
$ws.Cells.Item(73, 13).Value = -4942.1763  # M73: -4776.6665 -> -4942.1763
$ws.Cells.Item(73, 14).Value = -6804.3335  # N73: -6611.1 -> -6804.3335
$ws.Cells.Item(102, 8).Value = 5952.7666  # H102: 5666.968 -> 5952.7666
$ws.Cells.Item(102, 9).Value = 1387.8  # I102: 1356.7273 -> 1387.8
$ws.Cells.Item(102, 10).Value = 15082.7  # J102: 16203.111 -> 15082.7
$ws.Cells.Item(102, 11).Value = 1387.8  # K102: 1356.7273 -> 1387.8
$ws.Cells.Item(102, 12).Value = 15082.7  # L102: 16203.111 -> 15082.7
$ws.Cells.Item(102, 13).Value = 234.2  # M102: 265.2727 -> 234.2
$ws.Cells.Item(102, 14).Value = -18326.7  # N102: -19447.111 -> -18326.7
$ws.Cells.Item(117, 8).Value = 115999  # H117: 0 -> 115999
$ws.Cells.Item(117, 10).Value = 115999  # J117: 0 -> 115999
$ws.Cells.Item(117, 12).Value = 115999  # L117: 0 -> 115999
$ws.Cells.Item(117, 14).Value = -122883  # N117: None -> -122883

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 868.25  # H16: 845.48 -> 868.25
$ws.Cells.Item(16, 9).Value = 868.25  # I16: 845.48 -> 868.25
$ws.Cells.Item(16, 11).Value = 868.25  # K16: 845.48 -> 868.25
$ws.Cells.Item(16, 13).Value = -698.25  # M16: -675.48 -> -698.25
$ws.Cells.Item(88, 8).Value = 51745  # H88: 52495 -> 51745
$ws.Cells.Item(88, 9).Value = 28495  # I88: 29990 -> 28495
$ws.Cells.Item(88, 10).Value = 74995  # J88: 75000 -> 74995
$ws.Cells.Item(88, 11).Value = 28495  # K88: 29990 -> 28495
$ws.Cells.Item(88, 12).Value = 74995  # L88: 75000 -> 74995
$ws.Cells.Item(88, 13).Value = -28067  # M88: -29562 -> -28067
$ws.Cells.Item(88, 14).Value = -75851  # N88: -75856 -> -75851
$ws.Cells.Item(91, 8).Value = 51745  # H91: 52495 -> 51745
$ws.Cells.Item(91, 9).Value = 28495  # I91: 29990 -> 28495
$ws.Cells.Item(91, 10).Value = 74995  # J91: 75000 -> 74995
$ws.Cells.Item(91, 11).Value = 28495  # K91: 29990 -> 28495
$ws.Cells.Item(91, 12).Value = 74995  # L91: 75000 -> 74995
$ws.Cells.Item(91, 13).Value = -27013  # M91: -28508 -> -27013
$ws.Cells.Item(91, 14).Value = -77959  # N91: -77964 -> -77959
$ws.Cells.Item(110, 8).Value = 68322  # H110: 78322 -> 68322
$ws.Cells.Item(110, 10).Value = 68322  # J110: 78322 -> 68322
$ws.Cells.Item(110, 12).Value = 68322  # L110: 78322 -> 68322
$ws.Cells.Item(110, 14).Value = -76502  # N110: -86502 -> -76502
$ws.Cells.Item(115, 8).Value = 139995  # H115: 119000 -> 139995
$ws.Cells.Item(115, 10).Value = 139995  # J115: 119000 -> 139995
$ws.Cells.Item(115, 12).Value = 139995  # L115: 119000 -> 139995
$ws.Cells.Item(115, 14).Value = -142345  # N115: -121350 -> -142345
$ws.Cells.Item(136, 8).Value = 55056.145  # H136: 55156.594 -> 55056.145
$ws.Cells.Item(136, 9).Value = 11843.333  # I136: 9818.299999999999 -> 11843.333
$ws.Cells.Item(136, 11).Value = 35529.999  # K136: 29454.9 -> 35529.999
$ws.Cells.Item(136, 13).Value = -32979.999  # M136: -26904.9 -> -32979.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 5000  # H18: 4250 -> 5000
$ws.Cells.Item(18, 9).Value = 0  # I18: 3500 -> 0
$ws.Cells.Item(18, 11).Value = 0  # K18: 3500 -> 0
$ws.Cells.Item(18, 13).ClearContents()  # M18: -3327 -> (cleared)
$ws.Cells.Item(26, 8).Value = 0  # H26: 19999 -> 0
$ws.Cells.Item(26, 9).Value = 0  # I26: 19999 -> 0
$ws.Cells.Item(26, 11).Value = 0  # K26: 19999 -> 0
$ws.Cells.Item(26, 13).ClearContents()  # M26: -19706 -> (cleared)
$ws.Cells.Item(39, 8).Value = 30495  # H39: 20826.666 -> 30495
$ws.Cells.Item(39, 9).Value = 0  # I39: 30485 -> 0
$ws.Cells.Item(39, 10).Value = 30495  # J39: 15997.5 -> 30495
$ws.Cells.Item(39, 11).Value = 0  # K39: 30485 -> 0
$ws.Cells.Item(39, 12).Value = 30495  # L39: 15997.5 -> 30495
$ws.Cells.Item(39, 13).ClearContents()  # M39: -30072 -> (cleared)
$ws.Cells.Item(39, 14).Value = -31321  # N39: -16823.5 -> -31321
$ws.Cells.Item(43, 8).Value = 0  # H43: 59999 -> 0
$ws.Cells.Item(43, 9).Value = 0  # I43: 59999 -> 0
$ws.Cells.Item(43, 11).Value = 0  # K43: 59999 -> 0
$ws.Cells.Item(43, 13).ClearContents()  # M43: -59850 -> (cleared)
$ws.Cells.Item(49, 8).Value = 49998  # H49: 0 -> 49998
$ws.Cells.Item(49, 10).Value = 49998  # J49: 0 -> 49998
$ws.Cells.Item(49, 12).Value = 49998  # L49: 0 -> 49998
$ws.Cells.Item(49, 14).Value = -50458  # N49: None -> -50458
$ws.Cells.Item(123, 8).Value = 90000  # H123: 0 -> 90000
$ws.Cells.Item(123, 10).Value = 90000  # J123: 0 -> 90000
$ws.Cells.Item(123, 12).Value = 90000  # L123: 0 -> 90000
$ws.Cells.Item(123, 14).Value = -99800  # N123: None -> -99800
$ws.Cells.Item(125, 8).Value = 83417.336  # H125: 87441 -> 83417.336
$ws.Cells.Item(125, 10).Value = 83417.336  # J125: 87441 -> 83417.336
$ws.Cells.Item(125, 12).Value = 83417.336  # L125: 87441 -> 83417.336
$ws.Cells.Item(125, 14).Value = -93257.336  # N125: -97281 -> -93257.336
$ws.Cells.Item(136, 8).Value = 1686966.6  # H136: 1795597.6 -> 1686966.6
$ws.Cells.Item(136, 9).Value = 1673.3914  # I136: 1666.4546 -> 1673.3914
$ws.Cells.Item(136, 10).Value = 5563141  # J136: 6180763 -> 5563141
$ws.Cells.Item(136, 11).Value = 5020.174199999999  # K136: 4999.3638 -> 5020.174199999999
$ws.Cells.Item(136, 12).Value = 16689423  # L136: 18542289 -> 16689423
$ws.Cells.Item(136, 13).Value = -2470.174199999999  # M136: -2449.3638 -> -2470.174199999999
$ws.Cells.Item(136, 14).Value = -16694523  # N136: -18547389 -> -16694523

